# conf scheduling: Pinned talks only need to be set in the talks sheet.
# Room view shows solutions.
#
# 1) Rename the "* view" sheets to "*s view" / "Theme tracks view" so they
#    read as plural / consistent with the renamed "theme track" concept.
# 2) Rename "Theme" -> "Theme track" in the constraint name/description and
#    in the corresponding column headers.
# 3) The "Rooms" and "Speakers" sheets used to also bake the demo pinned
#    talk ("S06") into their own grids. That's redundant - the Talks sheet
#    is the single source of truth for pinning - so clear those leftover
#    indicators. The generated "*s view" sheets keep showing the solution.

$wb = $excel.ActiveWorkbook

# --- 1) Rename the view sheets ---------------------------------------
$wb.Worksheets.Item("Room view").Name = "Rooms view"
$wb.Worksheets.Item("Speaker view").Name = "Speakers view"
$wb.Worksheets.Item("Theme view").Name = "Theme tracks view"
$wb.Worksheets.Item("Sector view").Name = "Sectors view"
$wb.Worksheets.Item("Content view").Name = "Contents view"

# --- 2) Rename "Theme" wording to "Theme track" -----------------------
$configuration = $wb.Worksheets.Item("Configuration")
$configuration.Range("A4").Value = "Theme track conflict"
$configuration.Range("C4").Value = "Soft penalty per common theme track of 2 talks that have an overlapping timeslot"

$talks = $wb.Worksheets.Item("Talks")
$talks.Range("E1").Value = "Theme track tags"

$themeTracksView = $wb.Worksheets.Item("Theme tracks view")
$themeTracksView.Range("A2").Value = "Theme track tag"

# --- 3) Pinned talks only need to be set in the Talks sheet ------------
# Clear the leftover pinned-talk indicator (and its highlight colour) from
# the Rooms sheet.
$rooms = $wb.Worksheets.Item("Rooms")
$roomsPinnedCell = $rooms.Range("D3")
$roomsPinnedCell.ClearContents()
$roomsPinnedCell.ClearFormats()

# Clear the leftover pinned-talk indicator from the Speakers sheet, but
# keep the "unavailable" highlight colour in place.
$speakers = $wb.Worksheets.Item("Speakers")
$speakers.Range("K11").ClearContents()
